# Auto-applied market price data refresh across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2220.75
$ws.Range("J19").Value = 2221
$ws.Range("L19").Value = 2221
$ws.Range("N19").Value = -2571

$ws.Range("H32").Value = 2893.3
$ws.Range("J32").Value = 3022.5
$ws.Range("L32").Value = 3022.5
$ws.Range("N32").Value = -3674.5

$ws.Range("H33").Value = 276.72726
$ws.Range("I33").Value = 197.28572
$ws.Range("J33").Value = 415.75
$ws.Range("K33").Value = 197.28572
$ws.Range("L33").Value = 415.75
$ws.Range("M33").Value = 31.71428
$ws.Range("N33").Value = -873.75

$ws.Range("H38").Value = 2737.5715
$ws.Range("I38").Value = 1868.8334
$ws.Range("J38").Value = 7950
$ws.Range("K38").Value = 5606.5002
$ws.Range("L38").Value = 23850
$ws.Range("M38").Value = -5234.5002
$ws.Range("N38").Value = -24594

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = ""

$ws.Range("H58").Value = 3867.4285
$ws.Range("I58").Value = 55
$ws.Range("K58").Value = 165
$ws.Range("M58").Value = -15

$ws.Range("H62").Value = 2997.5
$ws.Range("I62").Value = 2997.5
$ws.Range("K62").Value = 2997.5
$ws.Range("M62").Value = -2373.5

$ws.Range("H65").Value = 2997.5
$ws.Range("I65").Value = 2997.5
$ws.Range("K65").Value = 14987.5
$ws.Range("M65").Value = -11867.5

$ws.Range("H86").Value = 12568.889
$ws.Range("I86").Value = 10758.333
$ws.Range("K86").Value = 10758.333
$ws.Range("M86").Value = -9635.333000000001

$ws.Range("H89").Value = 12568.889
$ws.Range("I89").Value = 10758.333
$ws.Range("K89").Value = 53791.665
$ws.Range("M89").Value = -48175.665

$ws.Range("H125").Value = 2822.6
$ws.Range("J125").Value = 2325.6667
$ws.Range("L125").Value = 20931.0003
$ws.Range("N125").Value = -25851.0003

$ws.Range("H132").Value = 2135.6667
$ws.Range("I132").Value = 2135.6667
$ws.Range("K132").Value = 6407.000100000001
$ws.Range("M132").Value = -3877.000100000001

$ws.Range("H138").Value = 3940.628
$ws.Range("J138").Value = 4312.75
$ws.Range("L138").Value = 12938.25
$ws.Range("N138").Value = -23218.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 976.5
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = ""

$ws.Range("H77").Value = 976.5
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = ""

$ws.Range("H132").Value = 2189.3635
$ws.Range("I132").Value = 2184.1428
$ws.Range("K132").Value = 6552.428400000001
$ws.Range("M132").Value = -4022.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2501.8333
$ws.Range("I105").Value = 2320
$ws.Range("K105").Value = 2320
$ws.Range("M105").Value = -573

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 454.875
$ws.Range("I5").Value = 397.25
$ws.Range("J5").Value = 512.5
$ws.Range("K5").Value = 397.25
$ws.Range("L5").Value = 512.5
$ws.Range("M5").Value = -285.25
$ws.Range("N5").Value = -736.5

$ws.Range("H7").Value = 149.84616
$ws.Range("I7").Value = 79
$ws.Range("K7").Value = 79
$ws.Range("M7").Value = 34

$ws.Range("H25").Value = 1950
$ws.Range("I25").Value = 1950
$ws.Range("K25").Value = 1950
$ws.Range("M25").Value = -1776

$ws.Range("H31").Value = 1333.1904
$ws.Range("I31").Value = 1665.3334
$ws.Range("J31").Value = 1277.8334
$ws.Range("K31").Value = 1665.3334
$ws.Range("L31").Value = 1277.8334
$ws.Range("M31").Value = -1370.3334
$ws.Range("N31").Value = -1867.8334

$ws.Range("H34").Value = 1333.1904
$ws.Range("I34").Value = 1665.3334
$ws.Range("J34").Value = 1277.8334
$ws.Range("K34").Value = 1665.3334
$ws.Range("L34").Value = 1277.8334
$ws.Range("M34").Value = -1463.3334
$ws.Range("N34").Value = -1681.8334

$ws.Range("H39").Value = 26509.9
$ws.Range("I39").Value = 11700.667
$ws.Range("J39").Value = 32856.715
$ws.Range("K39").Value = 11700.667
$ws.Range("L39").Value = 32856.715
$ws.Range("M39").Value = -11309.667
$ws.Range("N39").Value = -33638.715

$ws.Range("H48").Value = 15000
$ws.Range("J48").Value = 15000
$ws.Range("L48").Value = 15000
$ws.Range("N48").Value = -15952

$ws.Range("H49").Value = 26509.9
$ws.Range("I49").Value = 11700.667
$ws.Range("J49").Value = 32856.715
$ws.Range("K49").Value = 11700.667
$ws.Range("L49").Value = 32856.715
$ws.Range("M49").Value = -11518.667
$ws.Range("N49").Value = -33220.715

$ws.Range("H58").Value = 4283.4287
$ws.Range("I58").Value = 2886
$ws.Range("K58").Value = 2886
$ws.Range("M58").Value = -2683

$ws.Range("H107").Value = 1339.2222
$ws.Range("J107").Value = 1888
$ws.Range("L107").Value = 1888
$ws.Range("N107").Value = -5728

$ws.Range("H133").Value = 62326
$ws.Range("J133").Value = 62326
$ws.Range("L133").Value = 62326
$ws.Range("N133").Value = -67386

$ws.Range("H136").Value = 4283.4287
$ws.Range("I136").Value = 2886
$ws.Range("K136").Value = 8658
$ws.Range("M136").Value = -6108

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 51
$ws.Range("J12").Value = 57
$ws.Range("L12").Value = 171
$ws.Range("N12").Value = -517

$ws.Range("H34").Value = 286.25
$ws.Range("J34").Value = 315
$ws.Range("L34").Value = 945
$ws.Range("N34").Value = -1113

$ws.Range("H55").Value = 575.25
$ws.Range("I55").Value = 301
$ws.Range("J55").Value = 666.6667
$ws.Range("K55").Value = 903
$ws.Range("L55").Value = 2000.0001
$ws.Range("M55").Value = -726
$ws.Range("N55").Value = -2354.0001

$ws.Range("H139").Value = 3332.5
$ws.Range("I139").Value = 3035.9
$ws.Range("K139").Value = 9107.700000000001
$ws.Range("M139").Value = -3967.700000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 95.35294
$ws.Range("I2").Value = 96.545456
$ws.Range("J2").Value = 93.166664
$ws.Range("K2").Value = 96.545456
$ws.Range("L2").Value = 93.166664
$ws.Range("M2").Value = 16.454544
$ws.Range("N2").Value = -319.166664

$ws.Range("H43").Value = 12000
$ws.Range("J43").Value = 12000
$ws.Range("L43").Value = 12000
$ws.Range("N43").Value = -12302

$ws.Range("H46").Value = 88750
$ws.Range("I46").Value = 98333.336
$ws.Range("K46").Value = 98333.336
$ws.Range("M46").Value = -98177.336

$ws.Range("H80").Value = 5016
$ws.Range("I80").Value = 3650
$ws.Range("J80").Value = 5699
$ws.Range("K80").Value = 3650
$ws.Range("L80").Value = 5699
$ws.Range("M80").Value = -2652
$ws.Range("N80").Value = -7695

$ws.Range("H83").Value = 5016
$ws.Range("I83").Value = 3650
$ws.Range("J83").Value = 5699
$ws.Range("K83").Value = 18250
$ws.Range("L83").Value = 28495
$ws.Range("M83").Value = -13258
$ws.Range("N83").Value = -38479

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = ""

$ws.Range("H102").Value = 2619
$ws.Range("I102").Value = 2212.625
$ws.Range("J102").Value = 3702.6667
$ws.Range("K102").Value = 2212.625
$ws.Range("L102").Value = 3702.6667
$ws.Range("M102").Value = -590.625
$ws.Range("N102").Value = -6946.6667

$ws.Range("H122").Value = 3199
$ws.Range("I122").Value = 3699
$ws.Range("K122").Value = 11097
$ws.Range("M122").Value = -8647

$ws.Range("H126").Value = 4715.8184
$ws.Range("I126").Value = 3456.5
$ws.Range("J126").Value = 4995.6665
$ws.Range("K126").Value = 10369.5
$ws.Range("L126").Value = 14986.9995
$ws.Range("M126").Value = -7899.5
$ws.Range("N126").Value = -19926.9995

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").Value = ""

$ws.Range("H63").Value = 90077
$ws.Range("I63").Value = 90077
$ws.Range("K63").Value = 90077
$ws.Range("M63").Value = -89328

$ws.Range("H66").Value = 90077
$ws.Range("I66").Value = 90077
$ws.Range("K66").Value = 270231
$ws.Range("M66").Value = -266487

$ws.Range("H122").Value = 5625.814
$ws.Range("I122").Value = 4216.143
$ws.Range("J122").Value = 6971.409
$ws.Range("K122").Value = 12648.429
$ws.Range("L122").Value = 20914.227
$ws.Range("M122").Value = -10198.429
$ws.Range("N122").Value = -25814.227

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = ""

$ws.Range("H100").Value = 50001450
$ws.Range("I100").Value = 50001450
$ws.Range("K100").Value = 100002900
$ws.Range("M100").Value = -100002359

$ws.Range("H126").Value = 2150
$ws.Range("I126").Value = 2150
$ws.Range("K126").Value = 6450
$ws.Range("M126").Value = -3980

$ws.Range("H136").Value = 3879.8823
$ws.Range("I136").Value = 3980.5833
$ws.Range("K136").Value = 11941.7499
$ws.Range("M136").Value = -9391.749899999999
